$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''59.453.73'
$ws.Range('E2').Value = '  +1.59%  '
$ws.Range('D3').Value = '''2.610.08'
$ws.Range('E3').Value = '  +1.47%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''537.53'
$ws.Range('E5').Value = '  +4.17%  '
$ws.Range('D6').Value = '''141.05'
$ws.Range('E6').Value = '  +2.07%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').Value = '''0.569'
$ws.Range('E8').Value = '  +1.24%  '
$ws.Range('D9').Value = '''2.618.64'
$ws.Range('E9').Value = '  +1.15%  '
$ws.Range('D10').Value = '''6.47'
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('D11').Value = '''0.103'
$ws.Range('E11').Value = '  +3.51%  '
$ws.Range('E12').Value = '  +2.24%  '
$ws.Range('E13').Value = '  +1.92%  '
$ws.Range('D14').Value = '''3.070.71'
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('D15').Value = '''59.375.67'
$ws.Range('E15').Value = '  +1.52%  '
$ws.Range('D16').Value = '''20.62'
$ws.Range('E16').Value = '  +1.86%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '''2.665.89'
$ws.Range('E17').Value = '  +3.65%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '''0.0000134'
$ws.Range('E18').Value = '  +2.13%  '
$ws.Range('D19').Value = '''345.11'
$ws.Range('E19').Value = '  +2.59%  '
$ws.Range('E20').Value = '  +1.68%  '
$ws.Range('D21').Value = '''10.15'
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('D22').Value = '''6.40'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').Value = '''67.12'
$ws.Range('E24').Value = '  +2.00%  '
$ws.Range('E25').Value = '  +0.80%  '
$ws.Range('D26').Value = '''0.409'
$ws.Range('E26').Value = '  +1.98%  '
$ws.Range('D27').Value = '''0.999'
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('D28').Value = '''7.24'
$ws.Range('E28').Value = '  +3.56%  '
$ws.Range('D29').Value = '''0.0₃0750'
$ws.Range('E29').Value = '  +5.87%  '
$ws.Range('E30').Value = '  +0.14%  '
$ws.Range('E31').Value = '  +5.98%  '
$ws.Range('D32').Value = '''5.87'
$ws.Range('E32').Value = '  -0.68%  '
$ws.Range('D33').Value = '''18.90'
$ws.Range('E33').Value = '  +1.44%  '
$ws.Range('E34').Value = '  +0.69%  '
$ws.Range('E35').Value = '  +2.36%  '
$ws.Range('D36').Value = '''1.12'
$ws.Range('E36').Value = '  +0.49%  '
$ws.Range('D37').Value = '''37.02'
$ws.Range('E37').Value = '  +2.39%  '
$ws.Range('D38').Value = '''0.843'
$ws.Range('E38').Value = '  +1.04%  '
$ws.Range('D39').Value = '''1.46'
$ws.Range('E39').Value = '  +2.46%  '
$ws.Range('D40').Value = '''0.840'
$ws.Range('E40').Value = '  +2.64%  '
$ws.Range('D41').Value = '''3.56'
$ws.Range('E41').Value = '  +1.62%  '
$ws.Range('D42').Value = '''277.32'
$ws.Range('E42').Value = '  +2.24%  '
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('D44').Value = '''0.600'
$ws.Range('E44').Value = '  +2.56%  '
$ws.Range('E45').Value = '  +1.87%  '
$ws.Range('D46').Value = '''10.74'
$ws.Range('E46').Value = '  +0.18%  '
$ws.Range('E47').Value = '  +2.17%  '
$ws.Range('D48').Value = '''1.953.56'
$ws.Range('E48').Value = '  -0.70%  '
$ws.Range('D50').Value = '''18.39'
$ws.Range('E50').Value = '  +4.43%  '
$ws.Range('D51').Value = '''4.52'
$ws.Range('E51').Value = '  +2.30%  '
